# Update countries & provincias Spain
# Refresh COVID-19 country stats and re-sort two countries that moved
# ahead of their neighbours in the source ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Last updated" banner (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 20 de Abril de 2020 a las 19:22"

# Estados Unidos (row 4) - refreshed totals
$ws.Cells.Item(4, 2).Value = 771197
$ws.Cells.Item(4, 3).Value = 6561
$ws.Cells.Item(4, 5).Value = 658352
$ws.Cells.Item(4, 7).Value = 781
$ws.Cells.Item(4, 8).Value = 41356

# Paises Bajos (row 17) - casos criticos revised
$ws.Cells.Item(17, 6).Value = 1158

# Africa block: Egipto now ranks ahead of Sudafrica (rows 53-54 swap)
$ws.Cells.Item(53, 1).Value = "Egipto"
$ws.Cells.Item(53, 2).Value = 3333
$ws.Cells.Item(53, 3).Value = 189
$ws.Cells.Item(53, 4).Value = 821
$ws.Cells.Item(53, 5).Value = 2262
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 11
$ws.Cells.Item(53, 8).Value = 250

$ws.Cells.Item(54, 1).Value = "Sudafrica"
$ws.Cells.Item(54, 2).Value = 3158
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 903
$ws.Cells.Item(54, 5).Value = 2201
$ws.Cells.Item(54, 6).Value = 36
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 54

# Tanzania now ranks ahead of Mali..Islas Feroe (rows 122-127 shift down)
$ws.Cells.Item(122, 1).Value = "Tanzania"
$ws.Cells.Item(122, 2).Value = 254
$ws.Cells.Item(122, 3).Value = 84
$ws.Cells.Item(122, 4).Value = 11
$ws.Cells.Item(122, 5).Value = 233
$ws.Cells.Item(122, 6).Value = 4
$ws.Cells.Item(122, 7).Value = 3
$ws.Cells.Item(122, 8).Value = 10

$ws.Cells.Item(123, 1).Value = "Mali"
$ws.Cells.Item(123, 2).Value = 224
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(123, 4).Value = 42
$ws.Cells.Item(123, 5).Value = 168
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 14

$ws.Cells.Item(124, 1).Value = "El Salvador"
$ws.Cells.Item(124, 2).Value = 218
$ws.Cells.Item(124, 3).Value = 17
$ws.Cells.Item(124, 4).Value = 46
$ws.Cells.Item(124, 5).Value = 165
$ws.Cells.Item(124, 6).Value = 2
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 7

$ws.Cells.Item(125, 1).Value = "Paraguay"
$ws.Cells.Item(125, 2).Value = 208
$ws.Cells.Item(125, 3).Value = 2
$ws.Cells.Item(125, 4).Value = 46
$ws.Cells.Item(125, 5).Value = 154
$ws.Cells.Item(125, 6).Value = 1
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 8

$ws.Cells.Item(126, 1).Value = "Jamaica"
$ws.Cells.Item(126, 2).Value = 196
$ws.Cells.Item(126, 3).Value = 23
$ws.Cells.Item(126, 4).Value = 27
$ws.Cells.Item(126, 5).Value = 164
$ws.Cells.Item(126, 6).Value = 0
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 5

$ws.Cells.Item(127, 1).Value = "Islas Feroe"
$ws.Cells.Item(127, 2).Value = 185
$ws.Cells.Item(127, 3).Value = 0
$ws.Cells.Item(127, 4).Value = 176
$ws.Cells.Item(127, 5).Value = 9
$ws.Cells.Item(127, 6).Value = 0
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 0

# Togo (row 147) - activos/recuperados/muertes revised
$ws.Cells.Item(147, 4).Value = 53
$ws.Cells.Item(147, 5).Value = 25
$ws.Cells.Item(147, 7).Value = 1
$ws.Cells.Item(147, 8).Value = 6

# Santa Lucia (row 190) - activos/recuperados revised
$ws.Cells.Item(190, 4).Value = 13
$ws.Cells.Item(190, 5).Value = 2
